$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value  = 162
$ws1.Range("F11").Value = 6082
$ws1.Range("F12").Value = 59
$ws1.Range("F14").Value = 502
$ws1.Range("F17").Value = 364
$ws1.Range("F22").Value = 156
$ws1.Range("F25").Value = 1024
$ws1.Range("F27").Value = 1839
$ws1.Range("F28").Value = 498

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 272

$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 260

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 260
$ws4.Range("F10").Value = 162
$ws4.Range("F13").Value = 6082
$ws4.Range("F14").Value = 59
$ws4.Range("F17").Value = 502
$ws4.Range("F20").Value = 364
$ws4.Range("F25").Value = 272
$ws4.Range("F32").Value = 156
$ws4.Range("F35").Value = 1024
$ws4.Range("F37").Value = 1839
$ws4.Range("F38").Value = 498
